# Remove the final slide of the deck ("Thank you! / Survey and feedback"),
# along with its associated notes page, relationships and Content_Types
# registration. This mirrors deleting the last slide (slide 21, id 350)
# from the Slides pane in PowerPoint.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)
$s.Delete()
